# Auto-generated edit script: refresh market-price derived columns (H:N)
# across the per-job sheets, per the scheduled-runner data update.
$wb = $excel.ActiveWorkbook

# ======================== Sheet: ALC ========================
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1660.25
$ws.Range("I12").Value = 1880
$ws.Range("K12").Value = 1880
$ws.Range("M12").Value = -1710
# Row 19
$ws.Range("H19").Value = 781.4545000000001
$ws.Range("I19").Value = 685.25
$ws.Range("J19").Value = 836.4286
$ws.Range("K19").Value = 685.25
$ws.Range("L19").Value = 836.4286
$ws.Range("M19").Value = -510.25
$ws.Range("N19").Value = -1186.4286
# Row 40
$ws.Range("H40").Value = 1733.8334
$ws.Range("I40").Value = 1414
$ws.Range("J40").Value = 3333
$ws.Range("K40").Value = 1414
$ws.Range("L40").Value = 3333
$ws.Range("M40").Value = -1239
$ws.Range("N40").Value = -3683
# Row 127
$ws.Range("H127").Value = 1255.0769
$ws.Range("I127").Value = 781.1429000000001
$ws.Range("K127").Value = 2343.4287
$ws.Range("M127").Value = 2616.5713

# ======================== Sheet: ARM ========================
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3511.8333
$ws.Range("I45").Value = 3294.5
$ws.Range("K45").Value = 3294.5
$ws.Range("M45").Value = -2917.5
# Row 61
$ws.Range("H61").Value = 2300.484
$ws.Range("I61").Value = 1264.8182
$ws.Range("K61").Value = 1264.8182
$ws.Range("M61").Value = -1052.8182
# Row 97
$ws.Range("H97").Value = 2322.8572
$ws.Range("I97").Value = 2226.6667
$ws.Range("K97").Value = 2226.6667
$ws.Range("M97").Value = -1730.6667
# Row 110
$ws.Range("H110").Value = 383.33334
$ws.Range("I110").Value = 383.33334
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 383.33334
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1661.66666
$ws.Range("N110").ClearContents()
# Row 136
$ws.Range("H136").Value = 2300.484
$ws.Range("I136").Value = 1264.8182
$ws.Range("K136").Value = 3794.4546
$ws.Range("M136").Value = -1244.4546

# ======================== Sheet: BSM ========================
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 305.66666
$ws.Range("I22").Value = 305.66666
$ws.Range("K22").Value = 305.66666
$ws.Range("M22").Value = -132.66666
# Row 64
$ws.Range("H64").Value = 55556372
$ws.Range("J64").Value = 420
$ws.Range("L64").Value = 420
$ws.Range("N64").Value = -870
# Row 67
$ws.Range("H67").Value = 55556372
$ws.Range("J67").Value = 420
$ws.Range("L67").Value = 420
$ws.Range("N67").Value = -1980
# Row 99
$ws.Range("H99").Value = 2244.4666
$ws.Range("I99").Value = 2097.3333
$ws.Range("J99").Value = 2465.1667
$ws.Range("K99").Value = 2097.3333
$ws.Range("L99").Value = 2465.1667
$ws.Range("M99").Value = -599.3332999999998
$ws.Range("N99").Value = -5461.1667
# Row 100
$ws.Range("H100").Value = 34821.5
$ws.Range("J100").Value = 34821.5
$ws.Range("L100").Value = 34821.5
$ws.Range("N100").Value = -36985.5
# Row 105
$ws.Range("H105").Value = 2634595.8
$ws.Range("I105").Value = 3637.182
$ws.Range("J105").Value = 6252164
$ws.Range("K105").Value = 3637.182
$ws.Range("L105").Value = 6252164
$ws.Range("M105").Value = -1890.182
$ws.Range("N105").Value = -6255658
# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# ======================== Sheet: CRP ========================
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 800
$ws.Range("I16").Value = 800
$ws.Range("K16").Value = 800
$ws.Range("M16").Value = -513
# Row 31
$ws.Range("H31").Value = 3951.36
$ws.Range("J31").Value = 3942.9375
$ws.Range("L31").Value = 3942.9375
$ws.Range("N31").Value = -4532.9375
# Row 34
$ws.Range("H34").Value = 3951.36
$ws.Range("J34").Value = 3942.9375
$ws.Range("L34").Value = 3942.9375
$ws.Range("N34").Value = -4346.9375
# Row 58
$ws.Range("H58").Value = 18306.104
$ws.Range("I58").Value = 1076.5238
$ws.Range("J58").Value = 63533.75
$ws.Range("K58").Value = 1076.5238
$ws.Range("L58").Value = 63533.75
$ws.Range("M58").Value = -873.5237999999999
$ws.Range("N58").Value = -63939.75
# Row 74
$ws.Range("H74").Value = 28184.777
$ws.Range("J74").Value = 28184.777
$ws.Range("L74").Value = 28184.777
$ws.Range("N74").Value = -29932.777
# Row 77
$ws.Range("H77").Value = 28184.777
$ws.Range("J77").Value = 28184.777
$ws.Range("L77").Value = 84554.33099999999
$ws.Range("N77").Value = -93290.33099999999
# Row 96
$ws.Range("H96").Value = 4407.6
$ws.Range("J96").Value = 4407.6
$ws.Range("L96").Value = 4407.6
$ws.Range("N96").Value = -9899.6
# Row 113
$ws.Range("H113").Value = 800
$ws.Range("I113").Value = 800
$ws.Range("K113").Value = 800
$ws.Range("M113").Value = 1370
# Row 134
$ws.Range("H134").Value = 1286.5714
$ws.Range("I134").Value = 1001.5
$ws.Range("J134").Value = 1666.6666
$ws.Range("K134").Value = 3004.5
$ws.Range("L134").Value = 4999.9998
$ws.Range("M134").Value = -469.5
$ws.Range("N134").Value = -10069.9998
# Row 136
$ws.Range("H136").Value = 18306.104
$ws.Range("I136").Value = 1076.5238
$ws.Range("J136").Value = 63533.75
$ws.Range("K136").Value = 3229.5714
$ws.Range("L136").Value = 190601.25
$ws.Range("M136").Value = -679.5713999999998
$ws.Range("N136").Value = -195701.25
# Row 141
$ws.Range("H141").Value = 19807.428
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 21442
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 21442
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -31802

# ======================== Sheet: CUL ========================
$ws = $wb.Worksheets.Item("CUL")
# Row 107
$ws.Range("H107").Value = 100000
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
# Row 113
$ws.Range("H113").Value = 25533.625
$ws.Range("I113").Value = 100000.5
$ws.Range("J113").Value = 711.3333
$ws.Range("K113").Value = 300001.5
$ws.Range("L113").Value = 2133.9999
$ws.Range("M113").Value = -297831.5
$ws.Range("N113").Value = -6473.9999
# Row 131
$ws.Range("H131").Value = 114440.65
$ws.Range("I131").Value = 776.5
$ws.Range("J131").Value = 122757.54
$ws.Range("K131").Value = 2329.5
$ws.Range("L131").Value = 368272.62
$ws.Range("M131").Value = 2710.5
$ws.Range("N131").Value = -378352.62

# ======================== Sheet: LTW ========================
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2934.5715
$ws.Range("I40").Value = 1992
$ws.Range("J40").Value = 3641.5
$ws.Range("K40").Value = 1992
$ws.Range("L40").Value = 3641.5
$ws.Range("M40").Value = -1856
$ws.Range("N40").Value = -3913.5
# Row 104
$ws.Range("H104").Value = 24498.889
$ws.Range("J104").Value = 24498.889
$ws.Range("L104").Value = 24498.889
$ws.Range("N104").Value = -31486.889

# ======================== Sheet: WVR ========================
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 3031716.8
$ws.Range("I107").Value = 1120
$ws.Range("J107").Value = 4547015
$ws.Range("K107").Value = 3360
$ws.Range("L107").Value = 13641045
$ws.Range("M107").Value = -1440
$ws.Range("N107").Value = -13644885

